# csv_format_standard_20200825.xlsx - small updates, added crosswalk
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Guide")

# NOTE: cells are edited in D20 -> D19 -> D18 -> D17 -> D13 order (bottom-up)
# to match the shared-string allocation order the source workbook ended up
# with (new/changed strings are appended to the table in write order).

# --- D20: All geographic coordinates... (rich text, bold EPSG 4326) ---
$d20 = "All geographic coordinates must be provided in WGS84 decimal format (EPSG 4326). Latitude and longitude must be provided as separate variables (i.e., in an adjacent column). For geolocated records, each row must contain coordinates. In cases where the entire file consists of measurements collected at a single location, a pair of geographic coordinates must be documented elsewhere if not provided as variables. Spatial data using different standards can be provided as a separate variable (i.e., in an adjacent column) but only in addition to WGS84 decimal format."
$ws.Range("D20").Value = $d20
$ws.Range("D20").Characters(70, 9).Font.Bold = $true

# --- D19: For data with multiple timestamped records... (plain text, rewritten) ---
$ws.Range("D19").Value = "For data with multiple timestamped records or when applicable, the variable name should specify if the measurement is the start, stop, or midpoint value, or it shoild be documented elsewhere."

# --- D18: All dates and times... (rich text, bold RFC 3339) ---
$d18 = 'All dates and times must be reported in Coordinated Universal Time (UTC) and follow the ISO 8601 standard (RFC 3339). Note that the use of "Z" and "T" characters are unnecessary. All times must be preceded with a date. In cases where the entire file consists of temporal data collected at a single date and time, the date and time must be documented elsewhere if not provided as a variable. Temporal data using different standards can be provided as a separate variable (i.e., in an adjacent column) but only in addition to UTC format.'
$ws.Range("D18").Value = $d18
$ws.Range("D18").Characters(108, 8).Font.Bold = $true

# --- D17: Measurement uncertainty ... (plain text, appended clause) ---
$ws.Range("D17").Value = "Measurement uncertainty, limits of detection, data quality indicators, and other flags pertaining to individual values should be reported as a separate variable (i.e., in an adjacent column) but only in addition to the reported values. If a coding system is used to describe the flags, it must be documented elsewhere."

# --- D13: Unique variable names... (plain append, no distinct run formatting needed) ---
$ws.Range("D13").Value = 'Unique variable names must be used. No spaces. Letters, numbers, a hyphen "-" (ASCII Code 45), and an underscore "_" (ASCII Code 95) are preferred in variable names.'

# --- Row height changes ---
$ws.Rows.Item(13).RowHeight = 34
$ws.Rows.Item(17).RowHeight = 68
$ws.Rows.Item(19).RowHeight = 51
$ws.Rows.Item(20).RowHeight = 119

# --- Sheet view: drop frozen/scrolled topLeftCell, update selection ---
$ws.Range("D15").Select()

# --- Workbook window geometry ---
$excel.Width = 19860
$excel.Height = 15400
